$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "['https://www.facebook.com/naturesocietysingapore/', 'https://www.instagram.com/naturesocietysingapore/?hl=en', 'https://www.facebook.com/groups/naturesocietysingapore/']"
$ws.Range("K4").Value = "['https://www.facebook.com/zerowastesg/', 'https://sg.linkedin.com/company/zerowastesg']"
$ws.Range("K11").Value = "['https://sg.linkedin.com/company/singapore-environment-council', 'https://www.instagram.com/secsingapore/?hl=en', 'https://www.facebook.com/SingaporeEnvironmentCouncil/', 'https://www.youtube.com/@SECSingapore']"
$ws.Range("K13").Value = "['https://www.facebook.com/repairkopitiam/', 'https://www.instagram.com/repairkopitiam/?hl=en', 'https://www.facebook.com/repairkopitiam/photos/']"
$ws.Range("K16").Value = "['https://www.facebook.com/cicadatreeecoplace/']"
$ws.Range("K17").Value = "['https://www.facebook.com/EarthObservatoryOfSingapore/']"
$ws.Range("K18").Value = "['https://www.facebook.com/groups/greendrinkssingapore/', 'https://www.facebook.com/greendrinkssg/?locale=ps_AF', 'https://www.youtube.com/channel/UCq_vO3-P1ide5sjEJdQAkag']"
$ws.Range("K21").Value = "['https://www.facebook.com/kampung.senang/', 'https://sg.linkedin.com/company/kampungsenang']"
$ws.Range("K24").Value = "['https://sg.linkedin.com/company/guildasia', 'https://www.facebook.com/guildsg/', 'https://www.youtube.com/watch?v=wnbMXZ4zuBM', 'https://www.facebook.com/Lionsforge/videos/guild-ground-up-innovation-labs-for-development-%E3%82%AE%E3%83%AB%E3%83%89-%E7%A4%BE%E5%8C%BA%E5%88%9B%E6%96%B0%E4%BC%9A%E9%A6%86-is-here-at-the-impact/491474829116857/', 'https://medium.com/@groundupinnovation/about', 'https://www.facebook.com/guildasia/']"
$ws.Range("K26").Value = "['https://www.facebook.com/TrashHeroSingapore/', 'https://www.instagram.com/wearetrashherosingapore/', 'https://www.facebook.com/TrashHeroSingapore/events/', 'https://www.facebook.com/TrashHeroSingapore/videos/1940092789613462/']"
$ws.Range("K27").Value = "['https://www.instagram.com/smallchangelastingimpact/', 'https://www.facebook.com/SmallChangeLastingImpact/', 'https://www.linkedin.com/pulse/power-small-actions-creating-meaningful-impact-lasting', 'https://medium.com/@contact_28344/the-power-of-social-impact-how-small-actions-can-create-lasting-change-b96cd9a49e55']"
$ws.Range("K31").Value = "['https://www.facebook.com/TeamSeagrass-172603406103907/', 'https://www.instagram.com/teamseagrass/', 'https://twitter.com/teamseagrass', 'https://www.flickr.com/groups/1047086@N21/']"
$ws.Range("K36").Value = "['https://www.facebook.com/ACRESasia/', 'https://sg.linkedin.com/company/acressg', 'https://www.facebook.com/ACRESasia/photos/']"
$ws.Range("K38").Value = "['https://www.facebook.com/groups/sgfoodrescue/', 'https://www.facebook.com/foodrescuesingapore/', 'https://www.instagram.com/sgfoodrescue/?hl=en']"
$ws.Range("K40").Value = "['https://sg.linkedin.com/company/eco-sim', 'https://www.facebook.com/EcoSIMClub/']"
$ws.Range("K41").Value = "['https://sg.linkedin.com/company/elsa-nuslaw', 'https://www.instagram.com/nuslawelsa/', 'https://www.facebook.com/elsa.nuslaw/', 'https://www.facebook.com/elsaualberta/', 'https://www.instagram.com/uottawaelsa/?hl=en']"
$ws.Range("K42").Value = "['https://www.facebook.com/nusvege/', 'https://www.instagram.com/nusvege/?hl=en']"
